$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add columns M (Low Stock = L - H) and N (M - 2*H) for rows 2..62 ---
$ws.Range("M2").Formula = "=L2-H2"
$ws.Range("N2").Formula = "=M2-2*H2"
$ws.Range("M3:M62").Formula = "=L3-H3"
$ws.Range("N3:N62").Formula = "=M3-2*H3"

# --- Data correction: row 22's "want" (L22) goes from 4 to 10 ---
$ws.Range("L22").Value2 = 10

# --- Conditional formatting ---
# Existing rule highlights I2:I62 > 0 in red; keep it but it will be
# renumbered (priority 2, new dxf slot) once the new rule is inserted.
# New rule highlights N2:N62 < 0 in red (restock warning).

# Recreate the I2:I62 rule first so it keeps its place at the top of the
# worksheet's conditionalFormatting list.
$ws.Range("I2:I62").FormatConditions.Delete()
$fcI = $ws.Range("I2:I62").FormatConditions.Add(1, 5, "0")

# Add the new N2:N62 rule second.
$fcN = $ws.Range("N2:N62").FormatConditions.Add(1, 6, "0")

# Colour the new rule's dxf first (so it claims the lower dxf index), then
# the I rule's dxf (claims the higher dxf index) - matches target dxfId
# assignment (I->2, N->1).
$fcN.Interior.Color = 255
$fcI.Interior.Color = 255

# Final priority order: N (new) = 1, I (existing) = 2.
$fcN.Priority = 1
$fcI.Priority = 2

# --- View: selection moved to N24, no longer scrolled/frozen at A25 ---
$ws.Range("N24").Select()
